$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch labels for 8A and 8B so that they are in chronological order:
# rows 83-86 (previously "8A") become "8B"; rows 87-89 (previously "8B") become "8A".
$ws.Range("A83").Value = "8B"
$ws.Range("A84").Value = "8B"
$ws.Range("A85").Value = "8B"
$ws.Range("A86").Value = "8B"
$ws.Range("A87").Value = "8A"
$ws.Range("A88").Value = "8A"
$ws.Range("A89").Value = "8A"

# Scroll the window so row 72 is at the top and select A90, matching the updated view state.
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A90").Select()
